# Daily attendance processing - 2026-01-29 19:46:42
# For every row in column G whose "Recorded By" value is
# "dnasr281@gmail.com, System", swap the order of the two names to
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 7))
    $v = $cell.Value2
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
